# Updates cryptos price/volume figures to the latest scrape snapshot.
# Each row: D = Price (text, may contain multiple "." thousand separators
# or look like a decimal number -- always force text with a leading
# apostrophe so Excel does not silently coerce it to a Number/Date).
# E = Volume(1h) change, already padded with spaces so it stays text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "'26.930.21"
$ws.Range("E2").Value = "  +0.94%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "'1.845.04"
$ws.Range("E3").Value = "  +0.96%  "

# Row 4: TetherUSD
$ws.Range("D4").Value = "'1.011"
$ws.Range("E4").Value = "  +0.22%  "

# Row 5: BNB
$ws.Range("D5").Value = "'309.50"
$ws.Range("E5").Value = "  +0.51%  "

# Row 7: XRP
$ws.Range("D7").Value = "'0.4775"
$ws.Range("E7").Value = "  +2.65%  "

# Row 8: Cardano
$ws.Range("E8").Value = "  +1.65%  "

# Row 9: Dogecoin
$ws.Range("D9").Value = "'0.07226"
$ws.Range("E9").Value = "  +1.24%  "

# Row 10: Polygon
$ws.Range("D10").Value = "'0.9259"
$ws.Range("E10").Value = "  +2.36%  "

# Row 11: Solana
$ws.Range("D11").Value = "'19.66"
$ws.Range("E11").Value = "  +1.25%  "

# Row 12: TRON
$ws.Range("D12").Value = "'0.07699"
$ws.Range("E12").Value = "  -0.90%  "

# Row 13: WrappedEther
$ws.Range("D13").Value = "'1.868.59"
$ws.Range("E13").Value = "  +1.74%  "

# Row 14: Polkadot
$ws.Range("D14").Value = "'5.318"
$ws.Range("E14").Value = "  +0.95%  "

# Row 15: Chainlink
$ws.Range("D15").Value = "'6.404"
$ws.Range("E15").Value = "  +0.92%  "

# Row 16: Litecoin
$ws.Range("D16").Value = "'88.82"
$ws.Range("E16").Value = "  +1.29%  "

# Row 17: BinanceUSD
$ws.Range("D17").Value = "'1.012"

# Row 18: ShibaInu
$ws.Range("D18").Value = "'0.000008643"
$ws.Range("E18").Value = "  +0.85%  "

# Row 19: Dai
$ws.Range("D19").Value = "'1.009"
$ws.Range("E19").Value = "  +0.25%  "

# Row 20: WrappedBTC
$ws.Range("D20").Value = "'26.969.07"
$ws.Range("E20").Value = "  +0.91%  "

# Row 21: Avalanche
$ws.Range("D21").Value = "'14.53"
$ws.Range("E21").Value = "  +2.33%  "

# Row 22: Uniswap
$ws.Range("D22").Value = "'5.059"
$ws.Range("E22").Value = "  +0.86%  "

# Row 23: Cosmos
$ws.Range("D23").Value = "'10.65"
$ws.Range("E23").Value = "  +0.92%  "

# Row 24: Toncoin
$ws.Range("D24").Value = "'1.927"
$ws.Range("E24").Value = "  +0.00%  "

# Row 25: Monero
$ws.Range("D25").Value = "'152.40"
$ws.Range("E25").Value = "  -0.17%  "

# Row 26: EthereumClassic
$ws.Range("E26").Value = "  +1.28%  "

# Row 27: LidoDAOToken
$ws.Range("D27").Value = "'1.992"
$ws.Range("E27").Value = "  +0.97%  "

# Row 28: BitcoinCash
$ws.Range("D28").Value = "'114.13"
$ws.Range("E28").Value = "  +0.07%  "

# Row 29: InternetComputer(DFINITY)
$ws.Range("D29").Value = "'4.937"
$ws.Range("E29").Value = "  +2.27%  "

# Row 30: Stellar
$ws.Range("D30").Value = "'0.08875"
$ws.Range("E30").Value = "  +0.73%  "

# Row 31: HuobiToken
$ws.Range("E31").Value = "  +5.34%  "

# Row 32: ARBITRUM
$ws.Range("D32").Value = "'1.170"
$ws.Range("E32").Value = "  +2.50%  "

# Row 33: ImmutableX
$ws.Range("D33").Value = "'0.7438"
$ws.Range("E33").Value = "  +1.62%  "

# Row 34: Filecoin
$ws.Range("D34").Value = "'4.489"
$ws.Range("E34").Value = "  +0.98%  "

# Row 35: RenderToken
$ws.Range("D35").Value = "'2.707"
$ws.Range("E35").Value = "  -0.38%  "

# Row 36: TrustWalletToken
$ws.Range("D36").Value = "'1.123"
$ws.Range("E36").Value = "  +4.21%  "

# Row 37: VeChain
$ws.Range("D37").Value = "'0.01958"
$ws.Range("E37").Value = "  +1.73%  "

# Row 38: Hedera
$ws.Range("D38").Value = "'0.05264"
$ws.Range("E38").Value = "  +2.56%  "

# Row 39: MXToken
$ws.Range("D39").Value = "'2.978"
$ws.Range("E39").Value = "  +1.71%  "

# Row 40: TheSandbox
$ws.Range("D40").Value = "'0.5186"
$ws.Range("E40").Value = "  +2.50%  "

# Row 41: FraxShare
$ws.Range("D41").Value = "'6.989"
$ws.Range("E41").Value = "  +1.55%  "

# Row 42: Algorand
$ws.Range("D42").Value = "'0.1509"
$ws.Range("E42").Value = "  +0.77%  "

# Row 43: Aptos
$ws.Range("D43").Value = "'8.187"
$ws.Range("E43").Value = "  +1.99%  "

# Row 44: EnergySwap
$ws.Range("D44").Value = "'10.58"
$ws.Range("E44").Value = "  +5.62%  "

# Row 45: Decentraland
$ws.Range("D45").Value = "'0.4712"
$ws.Range("E45").Value = "  +0.95%  "

# Row 46: PaxDollar
$ws.Range("D46").Value = "'1.011"
$ws.Range("E46").Value = "  +0.28%  "

# Row 47: Quant
$ws.Range("D47").Value = "'101.53"
$ws.Range("E47").Value = "  +3.32%  "

# Row 48: NEARProtocol
$ws.Range("D48").Value = "'1.602"
$ws.Range("E48").Value = "  +2.45%  "

# Row 49: Aave
$ws.Range("D49").Value = "'66.03"
$ws.Range("E49").Value = "  +3.12%  "

# Row 50: Cronos
$ws.Range("D50").Value = "'0.06024"
$ws.Range("E50").Value = "  -0.24%  "

# Row 51: EOS
$ws.Range("D51").Value = "'0.8850"
$ws.Range("E51").Value = "  +3.76%  "
